# RMA Complete Flow (Issue Credit) - SO TO RMA Receipt To Create Credit Memo
#
# The underlying Provar automation run minted a fresh batch of RMA numbers
# (RMA-478W, RMA-K3US, RMA-VCUF, RMA-PU7Q). The "RMA Details Maintenance
# Grid" sheet is the active/"current" data grid for the flow, and it gets
# repointed at the newest group (RMA-PU7Q) - Sales Order Line, Shipper Line
# and the Salesforce record Id for each of the 3 detail rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (Pro-Stock Product Lot Track line)
$ws.Range("E2").Value = "RMA-PU7Q-001"
$ws.Range("F2").Value = "RMA-PU7Q-1-1"
$ws.Range("J2").Value = "a7s5f000000xMsgAAE"

# Row 3 (Pro-Stock Product - Stock Product No Track line)
$ws.Range("E3").Value = "RMA-PU7Q-002"
$ws.Range("F3").Value = "RMA-PU7Q-1-2"
$ws.Range("J3").Value = "a7s5f000000xMshAAE"

# Row 4 (Pro-Stock Product serial Track line)
$ws.Range("E4").Value = "RMA-PU7Q-003"
$ws.Range("F4").Value = "RMA-PU7Q-1-3"
$ws.Range("J4").Value = "a7s5f000000xMsiAAE"

# The Shipper Line (F) and Id (J) columns are best-fit; re-autofit them now
# that the displayed text has changed.
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null
